# 21st jan 2020 run
# Append 5 new data rows (rows 21-25) for the 2020-01-21 scrape run,
# mirroring the layout/formatting already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date(serial), city, type, size, price, unit_price, land_size, count
$newRows = @(
    @{ Row = 21; Date = 43851; City = "Algyo";      Type = "House";  Size = 116.2941176470588; Price = 32385294.11764706; UnitPrice = 296691.8471062569; LandSize = 468;               Count = 17 },
    @{ Row = 22; Date = 43851; City = "Morahalom";  Type = "House";  Size = 103.9;              Price = 23926000;          UnitPrice = 232894.34794284;   LandSize = 179.0333333333333; Count = 30 },
    @{ Row = 23; Date = 43851; City = "Szeged";     Type = "Flat";   Size = 65.48168892718655;  Price = 27897462.30073244; UnitPrice = 429922.6476606705; LandSize = $null;              Count = 2321 },
    @{ Row = 24; Date = 43851; City = "Szeged";     Type = "Garage"; Size = 17.88607594936709;  Price = 4128101.265822785; UnitPrice = 239682.3778796598; LandSize = $null;              Count = 79 },
    @{ Row = 25; Date = 43851; City = "Szeged";     Type = "House";  Size = 166.7172413793103;  Price = 58262137.93103448; UnitPrice = 1062076.207434499; LandSize = 458.5724137931035;  Count = 725 }
)

# Seed the formatting for column A on the new rows by copying the existing
# date cell's format (number format, font, border, alignment) in one shot
# so the new cells reuse the same cell style (s="2") instead of minting new
# near-duplicate styles.
$ws.Range("A20").Copy() | Out-Null
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.City
    $ws.Cells.Item($row, 3).Value = $r.Type
    $ws.Cells.Item($row, 4).Value = $r.Size
    $ws.Cells.Item($row, 5).Value = $r.Price
    $ws.Cells.Item($row, 6).Value = $r.UnitPrice
    if ($null -ne $r.LandSize) {
        $ws.Cells.Item($row, 7).Value = $r.LandSize
    }
    $ws.Cells.Item($row, 8).Value = $r.Count
}
